# "finish dev of item buff"
# - intInit!activeShrineStageEventTime (row 6, col B) changes from 7 -> 120
# - intInit!allianceRegionMapBaseTimePerGrid row (row 11) is removed entirely,
#   shifting the following rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intInit")
$ws.Activate()

# Bump activeShrineStageEventTime from 7 to 120
$ws.Range("B6").Value = 120

# Remove the obsolete allianceRegionMapBaseTimePerGrid row, shifting rows 12:19 up to 11:18
$ws.Rows.Item(11).Delete()

# Reflect the new selection state left behind by the row deletion
$ws.Range("A11:XFD11").Select()

# Best-effort: restore the window position recorded in the workbook view
$win = $excel.Windows.Item(1)
$win.Left = 15840
$win.Top = 2860
